$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.757.24'
$ws.Range("E2").Value = '  -0.12%  '

# Row 3
$ws.Range("D3").Value = '2.313.35'
$ws.Range("E3").Value = '  +0.34%  '

# Row 4
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.78'
$ws.Range("E5").Value = '  -1.97%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '106.99'
$ws.Range("E6").Value = '  +2.33%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.628'
$ws.Range("E7").Value = '  -0.14%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.607'
$ws.Range("E9").Value = '  +0.34%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.18'
$ws.Range("E10").Value = '  +0.67%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0914'
$ws.Range("E11").Value = '  +0.43%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.37'
$ws.Range("E12").Value = '  -1.82%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.991'
$ws.Range("E14").Value = '  -2.01%  '

# Row 15
$ws.Range("E15").Value = '  -0.65%  '

# Row 16
$ws.Range("D16").Value = '2.664.50'
$ws.Range("E16").Value = '  +0.32%  '

# Row 17
$ws.Range("D17").Value = '2.324.42'
$ws.Range("E17").Value = '  +0.79%  '

# Row 18
$ws.Range("D18").Value = '42.933.43'
$ws.Range("E18").Value = '  +0.45%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.48'
$ws.Range("E19").Value = '  -1.09%  '

# Row 20
$ws.Range("E20").Value = '  -0.35%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.08'
$ws.Range("E21").Value = '  -11.92%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.62'
$ws.Range("E22").Value = '  -0.49%  '

# Row 23
$ws.Range("E23").Value = '  -2.11%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '265.75'
$ws.Range("E24").Value = '  -0.50%  '

# Row 25
$ws.Range("E25").Value = '  +1.00%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.14%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.75'
$ws.Range("E27").Value = '  +14.25%  '

# Row 28
$ws.Range("E28").Value = '  +0.29%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.29'
$ws.Range("E29").Value = '  -2.44%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.56'
$ws.Range("E30").Value = '  +2.66%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.45'
$ws.Range("E31").Value = '  -0.95%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '166.30'
$ws.Range("E32").Value = '  +0.08%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0876'
$ws.Range("E33").Value = '  -0.98%  '

# Row 34
$ws.Range("E34").Value = '  +5.26%  '

# Row 35
$ws.Range("E35").Value = '  -1.23%  '

# Row 36
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.69'
$ws.Range("E36").Value = '  +2.47%  '

# Row 37
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.112'
$ws.Range("E37").Value = '  -1.52%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0358'
$ws.Range("E38").Value = '  +0.96%  '

# Row 39
$ws.Range("E39").Value = '  +4.83%  '

# Row 40
$ws.Range("E40").Value = '  -1.50%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.61'
$ws.Range("E41").Value = '  +0.69%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '103.67'
$ws.Range("E42").Value = '  +8.38%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '71.02'
$ws.Range("E43").Value = '  +0.65%  '

# Row 44
$ws.Range("E44").Value = '  +1.61%  '

# Row 45
$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.93'
$ws.Range("E45").Value = '  +5.15%  '

# Row 46
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.01'
$ws.Range("E46").Value = '  +0.15%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '112.43'
$ws.Range("E47").Value = '  -3.04%  '

# Row 48
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '1.654.36'
$ws.Range("E48").Value = '  -2.71%  '

# Row 49
$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '76.94'
$ws.Range("E49").Value = '  -5.53%  '

# Row 50
$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.24'
$ws.Range("E50").Value = '  +1.12%  '

# Row 51
$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.79'
$ws.Range("E51").Value = '  -0.39%  '

